# This script reproduces the commit "getting ready to add second function":
#  - A1 stops being the (broken) =MODE(...) formula and becomes the plain
#    text label "Time Step" (matching the other header cells, which are
#    already plain text).
#  - The per-row Sensor 1 / Sensed State simulation results in columns
#    B, C and D are refreshed with their new values.
#  - The custom column width formatting (<cols>) that used to cover
#    columns A:D is dropped, so the sheet goes back to default widths.
#
# To drop the sheet-level column width formatting cleanly (Excel has no
# direct "remove the <cols> block" call), we build the new sheet content
# on a brand-new worksheet -- which naturally starts out with no column
# customization -- and then swap it in for the old sheet under the same
# name, instead of editing the old sheet (and its column formatting) in
# place.

$wb = $excel.ActiveWorkbook
$origName = $wb.ActiveSheet.Name

$newWs = $wb.Worksheets.Add()

# Header row (unchanged text, only A1 changes from a formula to text)
$newWs.Range("A1").Value = "Time Step"
$newWs.Range("B1").Value = "Truth State"
$newWs.Range("C1").Value = "Sensor 1"
$newWs.Range("D1").Value = "Sensed State"

# Time Step, Truth State, Sensor 1, Sensed State for rows 2..102
$data = @(
    @(0,3,1,3),
    @(1,3,1,3),
    @(2,1,1,1),
    @(3,1,0,1),
    @(4,1,0,1),
    @(5,0,0,1),
    @(6,0,0,1),
    @(7,0,0,1),
    @(8,0,0,1),
    @(9,0,0,1),
    @(10,0,0,1),
    @(11,0,0,1),
    @(12,0,0,1),
    @(13,0,0,1),
    @(14,0,0,1),
    @(15,0,0,1),
    @(16,0,0,1),
    @(17,0,0,1),
    @(18,0,0,1),
    @(19,0,0,1),
    @(20,0,0,1),
    @(21,0,0,1),
    @(22,0,0,1),
    @(23,0,0,1),
    @(24,0,0,1),
    @(25,0,0,1),
    @(26,0,0,1),
    @(27,0,0,1),
    @(28,0,0,1),
    @(29,0,0,1),
    @(30,0,0,1),
    @(31,0,0,1),
    @(32,0,0,1),
    @(33,0,0,1),
    @(34,0,0,1),
    @(35,0,0,1),
    @(36,0,0,1),
    @(37,0,0,1),
    @(38,0,0,1),
    @(39,0,0,1),
    @(40,0,0,1),
    @(41,0,0,1),
    @(42,0,0,1),
    @(43,0,0,1),
    @(44,0,0,1),
    @(45,0,0,1),
    @(46,0,0,1),
    @(47,0,0,1),
    @(48,0,0,1),
    @(49,0,0,1),
    @(50,0,0,1),
    @(51,3,1,3),
    @(52,1,1,3),
    @(53,1,0,1),
    @(54,1,0,1),
    @(55,0,0,1),
    @(56,0,0,1),
    @(57,0,0,1),
    @(58,0,0,1),
    @(59,0,0,1),
    @(60,0,0,1),
    @(61,0,0,1),
    @(62,0,0,1),
    @(63,0,0,1),
    @(64,0,0,1),
    @(65,0,0,1),
    @(66,0,0,1),
    @(67,0,0,1),
    @(68,0,0,1),
    @(69,0,0,1),
    @(70,0,0,1),
    @(71,0,0,1),
    @(72,0,0,1),
    @(73,0,0,1),
    @(74,0,0,1),
    @(75,0,0,1),
    @(76,0,0,1),
    @(77,0,0,1),
    @(78,0,0,1),
    @(79,0,0,1),
    @(80,0,0,1),
    @(81,0,0,1),
    @(82,0,0,1),
    @(83,0,0,1),
    @(84,0,0,1),
    @(85,0,0,1),
    @(86,0,0,1),
    @(87,0,0,1),
    @(88,0,0,1),
    @(89,0,0,1),
    @(90,0,0,1),
    @(91,0,0,1),
    @(92,0,0,1),
    @(93,0,0,1),
    @(94,0,0,1),
    @(95,0,0,1),
    @(96,0,0,1),
    @(97,0,0,1),
    @(98,0,0,1),
    @(99,0,0,1),
    @(100,0,0,1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $newWs.Cells.Item($row, 1).Value = $vals[0]
    $newWs.Cells.Item($row, 2).Value = $vals[1]
    $newWs.Cells.Item($row, 3).Value = $vals[2]
    $newWs.Cells.Item($row, 4).Value = $vals[3]
}

# Swap the freshly built sheet in for the original one, keeping the name.
$orig = $wb.Worksheets.Item($origName)
[void]$orig.Delete()
$newWs.Name = $origName
